$d = $word.ActiveDocument

# Update the title/date line (unique text, safe to Find/Replace)
$d.Content.Find.Execute("2025-12-09 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-10 Wednesday", 2)

# Update the division problems in the table. Cell access (by row/column)
# is used instead of a global Find/Replace because several problems share
# identical text (e.g. two "479÷6=" and two "730÷3=" cells) and must be
# replaced with different values depending on position.
$t = $d.Tables.Item(1)

$updates = @{
    1  = @{ 1 = "813÷2="; 2 = "502÷7="; 3 = "423÷3="; 4 = "445÷5="; 5 = "490÷6=" }
    5  = @{ 1 = "364÷4="; 2 = "905÷6="; 3 = "305÷9="; 4 = "184÷6="; 5 = "812÷7=" }
    9  = @{ 1 = "472÷4="; 2 = "403÷5="; 3 = "457÷8="; 4 = "378÷9="; 5 = "966÷4=" }
    13 = @{ 1 = "970÷7="; 2 = "479÷7="; 3 = "823÷4="; 4 = "176÷3="; 5 = "452÷7=" }
    17 = @{ 1 = "126÷7="; 2 = "160÷5="; 3 = "335÷2="; 4 = "222÷2="; 5 = "221÷5=" }
}

foreach ($rowIndex in $updates.Keys) {
    $cols = $updates[$rowIndex]
    foreach ($colIndex in $cols.Keys) {
        $cell = $t.Cell($rowIndex, $colIndex)
        $newText = $cols[$colIndex]
        $cellRange = $cell.Range
        $cellRange.MoveEnd(1, -1) | Out-Null
        $cellRange.Text = $newText
    }
}
